$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; the existing "test" column (C) shifts to D
$ws.Columns("C").Insert()

# Set header for the newly inserted column C and give it an explicit width
$ws.Range("C1").Value = "raw_value"
$ws.Columns("C").ColumnWidth = 8.62

# Update the active selection to A3 (matches post-edit workbook state)
$ws.Range("A3").Select()
